$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.832.91'
$ws.Range('E2').Value = '  -0.40%  '

$ws.Range('D3').Value = '3.428.76'
$ws.Range('E3').Value = '  -0.29%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '410.83'
$ws.Range('E5').Value = '  +0.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.96'
$ws.Range('E6').Value = '  +0.87%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.635'
$ws.Range('E7').Value = '  +2.39%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.736'
$ws.Range('E9').Value = '  -2.46%  '

$ws.Range('E10').Value = '  -2.21%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '43.88'
$ws.Range('E11').Value = '  +1.50%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000226'
$ws.Range('E12').Value = '  +12.48%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.35'
$ws.Range('E13').Value = '  +7.19%  '

$ws.Range('D14').Value = '3.976.62'
$ws.Range('E14').Value = '  -0.37%  '

$ws.Range('E15').Value = '  +0.32%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.33'
$ws.Range('E16').Value = '  +4.93%  '

$ws.Range('D17').Value = '3.446.63'
$ws.Range('E17').Value = '  -0.63%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.45'
$ws.Range('E18').Value = '  +8.47%  '

$ws.Range('B19').Value = 'Polygon'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.09'
$ws.Range('E19').Value = '  +4.07%  '

$ws.Range('D20').Value = '61.908.77'
$ws.Range('E20').Value = '  -0.36%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '497.01'
$ws.Range('E21').Value = '  +34.24%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '93.17'
$ws.Range('E22').Value = '  +7.48%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.34'
$ws.Range('E23').Value = '  +4.83%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.45'
$ws.Range('E24').Value = '  +1.70%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.36'
$ws.Range('E25').Value = '  +5.52%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '35.03'
$ws.Range('E26').Value = '  +10.73%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.21'
$ws.Range('E27').Value = '  +10.45%  '

$ws.Range('E28').Value = '  -0.26%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.65'
$ws.Range('E29').Value = '  -0.79%  '

$ws.Range('E30').Value = '  +3.42%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.70'
$ws.Range('E31').Value = '  -0.30%  '

$ws.Range('E32').Value = '  -0.53%  '

$ws.Range('E33').Value = '  -2.36%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '42.03'
$ws.Range('E34').Value = '  -4.56%  '

$ws.Range('E35').Value = '  +0.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0502'
$ws.Range('E36').Value = '  +1.76%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '55.06'
$ws.Range('E37').Value = '  +5.84%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.49'
$ws.Range('E38').Value = '  +3.78%  '

$ws.Range('E39').Value = '  -0.06%  '

$ws.Range('E40').Value = '  +4.83%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.74'
$ws.Range('E41').Value = '  +17.59%  '

$ws.Range('E42').Value = '  +1.01%  '

$ws.Range('E43').Value = '  +1.45%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '145.41'
$ws.Range('E44').Value = '  +1.58%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.39'
$ws.Range('E45').Value = '  +10.34%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.09'
$ws.Range('E46').Value = '  +5.92%  '

$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.38'
$ws.Range('E47').Value = '  +22.58%  '

$ws.Range('B48').Value = 'Celestia'
$ws.Range('C48').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.72'
$ws.Range('E48').Value = '  +0.32%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.18'
$ws.Range('E49').Value = '  +7.32%  '

$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '120.64'
$ws.Range('E50').Value = '  +30.11%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.142'
$ws.Range('E51').Value = '  +17.21%  '
